$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped crypto data.
# Price cells must remain plain text (matching the source inlineStr cells),
# so we force a Text number format before assigning, then restore the default
# "Normal" style so no residual formatting is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.659.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.474.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.34%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.554'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.00%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.516'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.32%  '

$ws.Range("E10").Value = '  +14.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.88'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.111'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.854.89'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.500.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.790'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.607.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("E19").Value = '  +2.72%  '

$ws.Range("E20").Value = '  +0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.85%  '

$ws.Range("E22").Value = '  +0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.85%  '

$ws.Range("E25").Value = '  +0.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.90'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("E28").Value = '  +4.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.86'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.51'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.01%  '

$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0767'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.53%  '

$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.32%  '

$ws.Range("E37").Value = '  -1.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.116'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.53%  '

$ws.Range("E40").Value = '  -1.85%  '

$ws.Range("E41").Value = '  +1.55%  '

$ws.Range("E42").Value = '  -3.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.978.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.51%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("E46").Value = '  +0.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.708.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.87'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.53%  '

